# Refresh the "cryptos" symbol table (GitHub Actions job, 2023-01-11 05:38 UTC):
# updates Price/Volume(1h) figures for the listed coins and swaps the
# BOLO / CoinbaseStockToken rows (46-47), matching the new rank order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must be stored as literal text (not coerced to
# a number/percentage) even when it looks numeric, e.g. "275.68" or "0.62%".
# Forcing a Text number format before the assignment keeps the literal string,
# then resetting the style back to Normal avoids leaving stray formatting on
# the cell (matching the original workbook, where these cells carry no style).
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "275.68"
Set-TextValue $ws.Range("E2") "0.62%"
Set-TextValue $ws.Range("D3") "27.15"
Set-TextValue $ws.Range("E3") "1.63%"
Set-TextValue $ws.Range("D4") "4.854"
Set-TextValue $ws.Range("E4") "-0.22%"
Set-TextValue $ws.Range("D5") "0.06397"
Set-TextValue $ws.Range("E5") "1.25%"
Set-TextValue $ws.Range("D6") "6.949"
Set-TextValue $ws.Range("E6") "0.89%"
Set-TextValue $ws.Range("D7") "1.188"
Set-TextValue $ws.Range("E7") "-7.29%"
Set-TextValue $ws.Range("D8") "0.8751"
Set-TextValue $ws.Range("E8") "0.37%"
Set-TextValue $ws.Range("D9") "0.1512"
Set-TextValue $ws.Range("E9") "3.58%"
Set-TextValue $ws.Range("D10") "0.05068"
Set-TextValue $ws.Range("E10") "0.32%"
Set-TextValue $ws.Range("D11") "0.07521"
Set-TextValue $ws.Range("E11") "1.76%"
Set-TextValue $ws.Range("E12") "0.86%"
Set-TextValue $ws.Range("D13") "0.08977"
Set-TextValue $ws.Range("E13") "-0.81%"
Set-TextValue $ws.Range("D14") "0.001562"
Set-TextValue $ws.Range("E14") "-1.19%"
Set-TextValue $ws.Range("D15") "0.0006381"
Set-TextValue $ws.Range("E15") "1.00%"
Set-TextValue $ws.Range("D16") "0.006177"
Set-TextValue $ws.Range("E16") "4.98%"
Set-TextValue $ws.Range("D17") "3.476"
Set-TextValue $ws.Range("E17") "0.73%"
Set-TextValue $ws.Range("D18") "3.307"
Set-TextValue $ws.Range("E18") "-0.46%"
Set-TextValue $ws.Range("E19") "-1.40%"
Set-TextValue $ws.Range("D21") "0.1348"
Set-TextValue $ws.Range("E21") "1.72%"
Set-TextValue $ws.Range("D22") "3.905"
Set-TextValue $ws.Range("E22") "0.24%"
Set-TextValue $ws.Range("E23") "1.23%"
Set-TextValue $ws.Range("D25") "0.001173"
Set-TextValue $ws.Range("E25") "-0.34%"
Set-TextValue $ws.Range("D26") "0.003858"
Set-TextValue $ws.Range("E26") "-9.48%"
Set-TextValue $ws.Range("E27") "0.09%"
Set-TextValue $ws.Range("E28") "14.70%"
Set-TextValue $ws.Range("D40") "0.04157"
Set-TextValue $ws.Range("E40") "3.09%"
Set-TextValue $ws.Range("D41") "0.006808"
Set-TextValue $ws.Range("E41") "2.39%"
Set-TextValue $ws.Range("E42") "0.78%"
Set-TextValue $ws.Range("E43") "4.84%"
Set-TextValue $ws.Range("E44") "-2.67%"
Set-TextValue $ws.Range("D45") "0.00005192"
Set-TextValue $ws.Range("E45") "-2.09%"
$ws.Range("B46").Value = "CoinbaseStockToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
Set-TextValue $ws.Range("D46") "0.01850"
Set-TextValue $ws.Range("E46") "-7.46%"
$ws.Range("B47").Value = "BOLO"
$ws.Range("C47").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
Set-TextValue $ws.Range("D47") "1.687"
Set-TextValue $ws.Range("E47") "16.05%"
